$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($addr in @("D5","D6","D8","D10","D11","D12","D13","D14","D18","D20","D22","D23","D24","D25","D26","D27","D30","D31","D32","D33","D34","D35","D36","D37","D38","D45","D46","D47","D48","D50","D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '68.036.71'
$ws.Range("E2").Value = '  +2.44%  '

$ws.Range("D3").Value = '3.596.56'
$ws.Range("E3").Value = '  +0.85%  '

$ws.Range("D5").Value = '203.40'
$ws.Range("E5").Value = '  +9.02%  '

$ws.Range("D6").Value = '566.18'
$ws.Range("E6").Value = '  -2.87%  '

$ws.Range("D7").Value = '3.591.76'
$ws.Range("E7").Value = '  +0.83%  '

$ws.Range("D8").Value = '0.620'
$ws.Range("E8").Value = '  +1.06%  '

$ws.Range("E9").Value = '  +0.08%  '

$ws.Range("D10").Value = '0.670'
$ws.Range("E10").Value = '  +0.11%  '

$ws.Range("D11").Value = '60.10'
$ws.Range("E11").Value = '  +13.09%  '

$ws.Range("D12").Value = '0.151'
$ws.Range("E12").Value = '  +3.57%  '

$ws.Range("D13").Value = '0.0000286'
$ws.Range("E13").Value = '  +9.45%  '

$ws.Range("D14").Value = '9.99'
$ws.Range("E14").Value = '  +1.81%  '

$ws.Range("D15").Value = '4.158.30'
$ws.Range("E15").Value = '  +0.55%  '

$ws.Range("D16").Value = '3.582.71'
$ws.Range("E16").Value = '  +0.34%  '

$ws.Range("E17").Value = '  +0.73%  '

$ws.Range("D18").Value = '18.97'
$ws.Range("E18").Value = '  +3.37%  '

$ws.Range("D19").Value = '67.825.90'
$ws.Range("E19").Value = '  +2.23%  '

$ws.Range("D20").Value = '12.37'
$ws.Range("E20").Value = '  +1.33%  '

$ws.Range("E21").Value = '  +1.53%  '

$ws.Range("D22").Value = '401.77'
$ws.Range("E22").Value = '  +1.74%  '

$ws.Range("D23").Value = '12.77'
$ws.Range("E23").Value = '  +13.40%  '

$ws.Range("D24").Value = '4.16'
$ws.Range("E24").Value = '  -4.18%  '

$ws.Range("D25").Value = '85.10'
$ws.Range("E25").Value = '  -1.00%  '

$ws.Range("D26").Value = '2.91'
$ws.Range("E26").Value = '  +0.17%  '

$ws.Range("D27").Value = '12.56'
$ws.Range("E27").Value = '  +0.67%  '

$ws.Range("E28").Value = '  +10.02%  '

$ws.Range("E29").Value = '  +1.03%  '

$ws.Range("D30").Value = '8.24'
$ws.Range("E30").Value = '  +16.17%  '

$ws.Range("D31").Value = '9.33'
$ws.Range("E31").Value = '  +4.26%  '

$ws.Range("D32").Value = '31.53'
$ws.Range("E32").Value = '  +1.35%  '

$ws.Range("D33").Value = '670.57'
$ws.Range("E33").Value = '  +7.97%  '

$ws.Range("D34").Value = '12.19'
$ws.Range("E34").Value = '  +0.10%  '

$ws.Range("D35").Value = '0.114'
$ws.Range("E35").Value = '  +0.66%  '

$ws.Range("D36").Value = '63.59'
$ws.Range("E36").Value = '  +0.16%  '

$ws.Range("D37").Value = '42.19'
$ws.Range("E37").Value = '  +1.82%  '

$ws.Range("D38").Value = '0.423'
$ws.Range("E38").Value = '  +6.56%  '

$ws.Range("E39").Value = '  -0.08%  '

$ws.Range("D40").Value = '3.280.21'
$ws.Range("E40").Value = '  +9.20%  '

$ws.Range("D41").Value = '0.0₃0762'
$ws.Range("E41").Value = '  +0.01%  '

$ws.Range("E42").Value = '  +11.77%  '

$ws.Range("E43").Value = '  +2.76%  '

$ws.Range("E44").Value = '  +8.48%  '

$ws.Range("D45").Value = '0.998'
$ws.Range("E45").Value = '  -0.19%  '

$ws.Range("D46").Value = '2.99'
$ws.Range("E46").Value = '  +28.91%  '

$ws.Range("D47").Value = '0.0417'
$ws.Range("E47").Value = '  +2.07%  '

$ws.Range("D48").Value = '2.76'
$ws.Range("E48").Value = '  +11.52%  '

$ws.Range("E49").Value = '  +1.03%  '

$ws.Range("D50").Value = '0.131'
$ws.Range("E50").Value = '  +0.15%  '

$ws.Range("D51").Value = '8.80'
$ws.Range("E51").Value = '  +2.73%  '
